# Auto-generated edit script: updates Leve profit-calculation cells
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board price data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 158.42857
$ws.Range("I9").Value = 158.42857
$ws.Range("K9").Value = 158.42857
$ws.Range("M9").Value = 10.57142999999999
$ws.Range("H64").Value = 5066.6665
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 5066.6665
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -4916
$ws.Range("H70").Value = 7637.625
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 9350.166999999999
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 28050.501
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -28590.501
$ws.Range("H73").Value = 7637.625
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 9350.166999999999
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 28050.501
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -29922.501
$ws.Range("H80").Value = 278.53333
$ws.Range("I80").Value = 265.25
$ws.Range("J80").Value = 293.7143
$ws.Range("K80").Value = 795.75
$ws.Range("L80").Value = 881.1428999999999
$ws.Range("M80").Value = 202.25
$ws.Range("N80").Value = -2877.1429
$ws.Range("H83").Value = 278.53333
$ws.Range("I83").Value = 265.25
$ws.Range("J83").Value = 293.7143
$ws.Range("K83").Value = 2387.25
$ws.Range("L83").Value = 2643.4287
$ws.Range("M83").Value = 2604.75
$ws.Range("N83").Value = -12627.4287
$ws.Range("H113").Value = 3161.6667
$ws.Range("I113").Value = 3161.6667
$ws.Range("K113").Value = 3161.6667
$ws.Range("M113").Value = 92.33329999999978
$ws.Range("H125").Value = 3033.8333
$ws.Range("I125").Value = 1365.5
$ws.Range("K125").Value = 12289.5
$ws.Range("M125").Value = -9829.5
$ws.Range("H131").Value = 5017.2856
$ws.Range("I131").Value = 624.3333
$ws.Range("K131").Value = 1872.9999
$ws.Range("M131").Value = 3167.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 9944.111000000001
$ws.Range("I110").Value = 11837.4
$ws.Range("J110").Value = 7577.5
$ws.Range("K110").Value = 11837.4
$ws.Range("L110").Value = 7577.5
$ws.Range("M110").Value = -9792.4
$ws.Range("N110").Value = -11667.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2096
$ws.Range("I86").Value = 1413
$ws.Range("K86").Value = 1413
$ws.Range("M86").Value = -290
$ws.Range("H89").Value = 2096
$ws.Range("I89").Value = 1413
$ws.Range("K89").Value = 7065
$ws.Range("M89").Value = -1449

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3650.6086
$ws.Range("J58").Value = 4844.7856
$ws.Range("L58").Value = 4844.7856
$ws.Range("N58").Value = -5250.7856
$ws.Range("H81").Value = 23000
$ws.Range("J81").Value = 23000
$ws.Range("L81").Value = 23000
$ws.Range("N81").Value = -24996
$ws.Range("H84").Value = 23000
$ws.Range("J84").Value = 23000
$ws.Range("L84").Value = 69000
$ws.Range("N84").Value = -78984
$ws.Range("H99").Value = 13199.607
$ws.Range("I99").Value = 9937.691999999999
$ws.Range("J99").Value = 16026.6
$ws.Range("K99").Value = 9937.691999999999
$ws.Range("L99").Value = 16026.6
$ws.Range("M99").Value = -8439.691999999999
$ws.Range("N99").Value = -19022.6
$ws.Range("H107").Value = 697.3333
$ws.Range("I107").Value = 648.4
$ws.Range("J107").Value = 942
$ws.Range("K107").Value = 648.4
$ws.Range("L107").Value = 942
$ws.Range("M107").Value = 1271.6
$ws.Range("N107").Value = -4782
$ws.Range("H126").Value = 13199.607
$ws.Range("I126").Value = 9937.691999999999
$ws.Range("J126").Value = 16026.6
$ws.Range("K126").Value = 29813.076
$ws.Range("L126").Value = 48079.8
$ws.Range("M126").Value = -27343.076
$ws.Range("N126").Value = -53019.8
$ws.Range("H134").Value = 2399.5757
$ws.Range("I134").Value = 1812.8334
$ws.Range("J134").Value = 3964.2222
$ws.Range("K134").Value = 5438.5002
$ws.Range("L134").Value = 11892.6666
$ws.Range("M134").Value = -2903.5002
$ws.Range("N134").Value = -16962.6666
$ws.Range("H136").Value = 3650.6086
$ws.Range("J136").Value = 4844.7856
$ws.Range("L136").Value = 14534.3568
$ws.Range("N136").Value = -19634.3568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1532.091
$ws.Range("J33").Value = 934
$ws.Range("L33").Value = 5604
$ws.Range("N33").Value = -6170

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1580.7894
$ws.Range("I102").Value = 366.81818
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 366.81818
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = 1255.18182
$ws.Range("N102").Value = -6494
$ws.Range("H126").Value = 4476.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2751.5454
$ws.Range("I7").Value = 2740.7778
$ws.Range("K7").Value = 2740.7778
$ws.Range("M7").Value = -2628.7778
$ws.Range("H82").Value = 2947.077
$ws.Range("I82").Value = 3219.1052
$ws.Range("K82").Value = 3219.1052
$ws.Range("M82").Value = -2858.1052
$ws.Range("H85").Value = 2947.077
$ws.Range("I85").Value = 3219.1052
$ws.Range("K85").Value = 3219.1052
$ws.Range("M85").Value = -1971.1052
$ws.Range("H100").Value = 2099.875
$ws.Range("I100").Value = 1933
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 1933
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -1392
$ws.Range("N100").Value = -3282
$ws.Range("H126").Value = 2751.5454
$ws.Range("I126").Value = 2740.7778
$ws.Range("K126").Value = 8222.3334
$ws.Range("M126").Value = -5752.3334
$ws.Range("H136").Value = 3999.3333
$ws.Range("I136").Value = 3599.4
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 10798.2
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -8248.200000000001
$ws.Range("N136").Value = -23097

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1228
$ws.Range("I96").Value = 1056.8334
$ws.Range("J96").Value = 1399.1666
$ws.Range("K96").Value = 1056.8334
$ws.Range("L96").Value = 1399.1666
$ws.Range("M96").Value = 316.1666
$ws.Range("N96").Value = -4145.1666
